# Swap the "lcoe_class" (column P) values for paired cost-class rows
# on the "solar" and "wind" sheets, per the 2025-08-13 BGR model update.

$wb = $excel.ActiveWorkbook

$solar = $wb.Worksheets.Item("solar")
$solar.Range("P5").Value = 2
$solar.Range("P6").Value = 4

$wind = $wb.Worksheets.Item("wind")
$wind.Range("P13").Value = 4
$wind.Range("P14").Value = 5
$wind.Range("P16").Value = 1
$wind.Range("P17").Value = 3
$wind.Range("P27").Value = 1
$wind.Range("P28").Value = 2
